# Fix excel upload error
# - Remove the sample/example data row (joe biden / Joever@gmail.com / Female)
#   that was left in row 2, together with its mailto hyperlink on A2 (the
#   hyperlink was causing an upload bug).
# - Keep the formatting (style) of A2 but clear its value/hyperlink.
# - Update the instructions note in E3 with an extra warning line about not
#   turning the gmail address into a clickable link.
# - Grow row 3 so the longer note still fits, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink attached to A2 (mailto:Joever@gmail.com) entirely.
$ws.Hyperlinks.Delete()

# Clear the sample data row (A2:C2) but keep cell formatting/style.
$ws.Range("A2:C2").ClearContents()

# Update the note text in E3 to include the new warning line.
$ws.Range("E3").Value = "*DO NOT DELETE TITLE ROW" + [char]10 + "- Gender is either Male or Femail" + [char]10 + "- Email is either @gmail.com or @fpt.edu.vn" + [char]10 + "- DO NOT make gmail into a url/link, it will bug"

# Grow row 3 to fit the extra line of text.
$ws.Rows.Item(3).RowHeight = 90

# Move the active selection to F3, matching the saved selection state.
$ws.Range("F3").Select()
